$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column A (rows 3-10) with 1000, continuing the existing series in A1:A2
$ws.Range("A3").Value = 1000
$ws.Range("A4").Value = 1000
$ws.Range("A5").Value = 1000
$ws.Range("A6").Value = 1000
$ws.Range("A7").Value = 1000
$ws.Range("A8").Value = 1000
$ws.Range("A9").Value = 1000
$ws.Range("A10").Value = 1000

# Move/record the active selection to A10, as it was left after entering values
$ws.Range("A10").Select()
